$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "Brasil", "01/01/2014", 0.001474320658766841),
    @(3, "Brasil", "01/01/2015", 0.2166497206675002),
    @(4, "Brasil", "01/01/2016", 0.1997052926159623),
    @(5, "Brasil", "01/01/2017", 0.2000173550022124),
    @(6, "Brasil", "01/01/2018", 0.2218892578236288),
    @(7, "Brasil", "01/01/2019", 0.2849038916940809),
    @(8, "Brasil", "01/01/2020", 0.3176212688798147),
    @(9, "Brasil", "01/01/2021", 0.389362858101849),
    @(10, "Brasil", "01/01/2022", 0.4803159880483661),
    @(11, "Brasil", "01/01/2023", 0.7045059749555662),
    @(12, "Brasil", "01/01/2024", 0.7406584305647893),
    @(13, "Nordeste", "01/01/2014", 0.001702829476764164),
    @(14, "Nordeste", "01/01/2015", 0.2786930553707959),
    @(15, "Nordeste", "01/01/2016", 0.2922254241262205),
    @(16, "Nordeste", "01/01/2017", 0.269560011872592),
    @(17, "Nordeste", "01/01/2018", 0.3086050738242937),
    @(18, "Nordeste", "01/01/2019", 0.3800713524743992),
    @(19, "Nordeste", "01/01/2020", 0.468542465255369),
    @(20, "Nordeste", "01/01/2021", 0.5294632190918449),
    @(21, "Nordeste", "01/01/2022", 0.6728854456429895),
    @(22, "Nordeste", "01/01/2023", 0.8967265550752545),
    @(23, "Nordeste", "01/01/2024", 1.008284356164113),
    @(24, "Sergipe", "01/01/2014", 0.0003164897488970332),
    @(25, "Sergipe", "01/01/2015", 0.1812486638720291),
    @(26, "Sergipe", "01/01/2016", 0.1731002802120558),
    @(27, "Sergipe", "01/01/2017", 0.2406581184562085),
    @(28, "Sergipe", "01/01/2018", 0.2878955832071318),
    @(29, "Sergipe", "01/01/2019", 0.4083701716603065),
    @(30, "Sergipe", "01/01/2020", 0.4713832671714768),
    @(31, "Sergipe", "01/01/2021", 0.5925616849552231),
    @(32, "Sergipe", "01/01/2022", 0.7949149714042992),
    @(33, "Sergipe", "01/01/2023", 1.115654925197335),
    @(34, "Sergipe", "01/01/2024", 1.197113711052504)
)

foreach ($row in $data) {
    $r = $row[0]
    $region = $row[1]
    $dateText = $row[2]
    $rate = $row[3]

    $ws.Cells.Item($r, 1).Value = $region

    $dc = $ws.Cells.Item($r, 2)
    $dc.Formula = '="' + $dateText + '"'
    $dc.Copy()
    $dc.PasteSpecial(-4163)

    $ws.Cells.Item($r, 3).Value = $rate
}
$ws.Application.CutCopyMode = $false

# Remove now-obsolete trailing rows (previously 35-43; table now ends at row 34)
for ($i = 43; $i -ge 35; $i--) {
    $ws.Rows.Item($i).EntireRow.Delete()
}
